$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33 (ALC)
$ws.Range("H33").Value = 3925.1428
$ws.Range("I33").Value = 3925.1428
$ws.Range("K33").Value = 3925.1428
$ws.Range("M33").Value = -3696.1428

# Row 42 (ALC)
$ws.Range("H42").Value = 48.166668
$ws.Range("I42").Value = 22.25
$ws.Range("K42").Value = 66.75
$ws.Range("M42").Value = 163.25

# Row 98 (ALC)
$ws.Range("H98").Value = 7148.5454
$ws.Range("I98").Value = 1980.6666
$ws.Range("K98").Value = 1980.6666
$ws.Range("M98").Value = -482.6666

# Row 101 (ALC)
$ws.Range("H101").Value = 1248.75
$ws.Range("I101").Value = 998.3333
$ws.Range("J101").Value = 2000
$ws.Range("K101").Value = 2994.9999
$ws.Range("L101").Value = 6000
$ws.Range("M101").Value = -1372.9999
$ws.Range("N101").Value = -9244

# Row 122 (ALC)
$ws.Range("H122").Value = 7148.5454
$ws.Range("I122").Value = 1980.6666
$ws.Range("K122").Value = 5941.9998
$ws.Range("M122").Value = -3491.9998

# Row 137 (ALC)
$ws.Range("H137").Value = 17497.25
$ws.Range("I137").Value = 14996.333
$ws.Range("K137").Value = 44988.999
$ws.Range("M137").Value = -42438.999

# Row 138 (ALC)
$ws.Range("H138").Value = 2210.9092
$ws.Range("I138").Value = 406.6
$ws.Range("J138").Value = 3714.5
$ws.Range("K138").Value = 1219.8
$ws.Range("L138").Value = 11143.5
$ws.Range("M138").Value = 3920.2
$ws.Range("N138").Value = -21423.5

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (ARM)
$ws.Range("H32").Value = 3583.1667
$ws.Range("I32").Value = 3583.1667
$ws.Range("K32").Value = 3583.1667
$ws.Range("M32").Value = -3296.1667

# Row 61 (ARM)
$ws.Range("H61").Value = 3066.3333
$ws.Range("I61").Value = 3066.3333
$ws.Range("K61").Value = 3066.3333
$ws.Range("M61").Value = -2854.3333

# Row 97 (ARM)
$ws.Range("H97").Value = 1005
$ws.Range("I97").Value = 1005
$ws.Range("K97").Value = 1005
$ws.Range("M97").Value = -509

# Row 102 (ARM)
$ws.Range("H102").Value = 4566
$ws.Range("I102").Value = 4879.2
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 4879.2
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = -3257.2
$ws.Range("N102").Value = -6244

# Row 122 (ARM)
$ws.Range("H122").Value = 2670.3333
$ws.Range("I122").Value = 2012
$ws.Range("J122").Value = 2999.5
$ws.Range("K122").Value = 6036
$ws.Range("L122").Value = 8998.5
$ws.Range("M122").Value = -3586
$ws.Range("N122").Value = -13898.5

# Row 132 (ARM)
$ws.Range("H132").Value = 1700.2858
$ws.Range("I132").Value = 1650.3334
$ws.Range("K132").Value = 4951.0002
$ws.Range("M132").Value = -2421.0002

# Row 136 (ARM)
$ws.Range("H136").Value = 3066.3333
$ws.Range("I136").Value = 3066.3333
$ws.Range("K136").Value = 9198.999899999999
$ws.Range("M136").Value = -6648.999899999999

$ws = $wb.Worksheets.Item("BSM")
# Row 105 (BSM)
$ws.Range("H105").Value = 9749.25
$ws.Range("I105").Value = 8999
$ws.Range("K105").Value = 8999
$ws.Range("M105").Value = -7252

# Row 134 (BSM)
$ws.Range("H134").Value = 5354.8
$ws.Range("I134").Value = 4721.5
$ws.Range("K134").Value = 14164.5
$ws.Range("M134").Value = -11629.5

$ws = $wb.Worksheets.Item("CRP")
# Row 99 (CRP)
$ws.Range("H99").Value = 12000
$ws.Range("I99").Value = 12000
$ws.Range("K99").Value = 12000
$ws.Range("M99").Value = -10502

# Row 100 (CRP)
$ws.Range("H100").Value = 99998.5
$ws.Range("J100").Value = 99998.5
$ws.Range("L100").Value = 99998.5
$ws.Range("N100").Value = -102162.5

# Row 108 (CRP)
$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("M108").ClearContents()
$ws.Range("N108").ClearContents()

# Row 126 (CRP)
$ws.Range("H126").Value = 12000
$ws.Range("I126").Value = 12000
$ws.Range("K126").Value = 36000
$ws.Range("M126").Value = -33530

# Row 132 (CRP)
$ws.Range("H132").Value = 2797.8
$ws.Range("I132").Value = 1663.3334
$ws.Range("K132").Value = 4990.0002
$ws.Range("M132").Value = -2460.0002

$ws = $wb.Worksheets.Item("CUL")
# Row 10 (CUL)
$ws.Range("H10").Value = 183.81818
$ws.Range("I10").Value = 72.333336
$ws.Range("K10").Value = 217.000008
$ws.Range("M10").Value = -78.00000800000001

# Row 16 (CUL)
$ws.Range("H16").Value = 833.3333
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 833.3333
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 2499.9999
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -2845.9999

$ws = $wb.Worksheets.Item("GSM")
# Row 122 (GSM)
$ws.Range("H122").Value = 3636.8
$ws.Range("I122").Value = 2035.6666
$ws.Range("J122").Value = 6038.5
$ws.Range("K122").Value = 6106.9998
$ws.Range("L122").Value = 18115.5
$ws.Range("M122").Value = -3656.9998
$ws.Range("N122").Value = -23015.5

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (LTW)
$ws.Range("H7").Value = 10916.167
$ws.Range("I7").Value = 11665.667
$ws.Range("K7").Value = 11665.667
$ws.Range("M7").Value = -11553.667

# Row 16 (LTW)
$ws.Range("H16").Value = 1075.5
$ws.Range("I16").Value = 934
$ws.Range("J16").Value = 1500
$ws.Range("K16").Value = 934
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = -764
$ws.Range("N16").Value = -1840

# Row 40 (LTW)
$ws.Range("H40").Value = 10789.833
$ws.Range("J40").Value = 10960.8
$ws.Range("L40").Value = 10960.8
$ws.Range("N40").Value = -11232.8

# Row 122 (LTW)
$ws.Range("H122").Value = 3580.6667
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

# Row 126 (LTW)
$ws.Range("H126").Value = 10916.167
$ws.Range("I126").Value = 11665.667
$ws.Range("K126").Value = 34997.001
$ws.Range("M126").Value = -32527.001

# Row 132 (LTW)
$ws.Range("H132").Value = 7601
$ws.Range("I132").Value = 7920.125
$ws.Range("K132").Value = 23760.375
$ws.Range("M132").Value = -21230.375

# Row 136 (LTW)
$ws.Range("H136").Value = 4417.8335
$ws.Range("I136").Value = 4601.4
$ws.Range("K136").Value = 13804.2
$ws.Range("M136").Value = -11254.2

$ws = $wb.Worksheets.Item("WVR")
# Row 107 (WVR)
$ws.Range("H107").Value = 1588
$ws.Range("I107").Value = 1449.6666
$ws.Range("K107").Value = 4348.9998
$ws.Range("M107").Value = -2428.9998

# Row 125 (WVR)
$ws.Range("H125").Value = 114999
$ws.Range("J125").Value = 114999
$ws.Range("L125").Value = 114999
$ws.Range("N125").Value = -124839

# Row 132 (WVR)
$ws.Range("H132").Value = 4261.6
$ws.Range("I132").Value = 3148.077
$ws.Range("K132").Value = 9444.231
$ws.Range("M132").Value = -6914.231

# Row 136 (WVR)
$ws.Range("H136").Value = 10360.667
$ws.Range("I136").Value = 11777.714
$ws.Range("J136").Value = 5401
$ws.Range("K136").Value = 35333.142
$ws.Range("L136").Value = 16203
$ws.Range("M136").Value = -32783.142
$ws.Range("N136").Value = -21303
